$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 7145989545
$ws.Range("E3").Value = 7145989545

$ws.Range("E7").Select()
